$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell C1 gets new label "H02"
$ws.Range("C1").Value = "H02"

# Fill in column C scores for rows 2-16 (Homework 02 column)
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 9.5
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 9.5
$ws.Range("C11").Value = 7
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = 5
$ws.Range("C14").Value = 7
$ws.Range("C15").Value = 10
$ws.Range("C16").Value = 8.5

# Row 10's extra-point columns (C:E) pick up row's main numeric style
# (s=4) instead of the distinct comma-decimal style (s=5) they had before
$ws.Range("C10:E10").NumberFormat = $ws.Range("B10").NumberFormat
$ws.Range("C10:E10").HorizontalAlignment = $ws.Range("B10").HorizontalAlignment

# Rows 11-16 previously had no cell at all in column C; typed values pick
# up the row's existing numeric style (s=4) matching column B
$ws.Range("C11:C16").NumberFormat = $ws.Range("B11").NumberFormat
$ws.Range("C11:C16").HorizontalAlignment = $ws.Range("B11").HorizontalAlignment

# Update B17 formula to include full row (COUNTA over B1:Q1)
$ws.Range("B17").Formula = "=COUNTA(B1:Q1) * 10"

# Update B22:B36 formulas to sum full row (Bx:Rx) instead of single cell
$ws.Range("B22").Formula = "=(SUM(B2:R2)/B17) * 40"
$ws.Range("B23").Formula = "=(SUM(B3:R3)/B17) * 40"
$ws.Range("B24").Formula = "=(SUM(B4:R4)/B17) * 40"
$ws.Range("B25").Formula = "=(SUM(B5:R5)/B17) * 40"
$ws.Range("B26").Formula = "=(SUM(B6:R6)/B17) * 40"
$ws.Range("B27").Formula = "=(SUM(B7:R7)/B17) * 40"
$ws.Range("B28").Formula = "=(SUM(B8:R8)/B17) * 40"
$ws.Range("B29").Formula = "=(SUM(B9:R9)/B17) * 40"
$ws.Range("B30").Formula = "=(SUM(B10:R10)/B17) * 40"
$ws.Range("B31").Formula = "=(SUM(B11:R11)/B17) * 40"
$ws.Range("B32").Formula = "=(SUM(B12:R12)/B17) * 40"
$ws.Range("B33").Formula = "=(SUM(B13:R13)/B17) * 40"
$ws.Range("B34").Formula = "=(SUM(B14:R14)/B17) * 40"
$ws.Range("B35").Formula = "=(SUM(B15:R15)/B17) * 40"
$ws.Range("B36").Formula = "=(SUM(B16:R16)/B17) * 40"

# Update selection to C8
$ws.Range("C8").Select()
